$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "CausalImpact는 어떻게 이벤트 효과를 추정할까?"
$ws.Range("E3").Value = "https://lumiamitie.github.io/data/causalimpact/"

$ws.Range("D4").Value = "Matplotlib 컬러명(color name), 팔레트(palette) 이름"
$ws.Range("E4").Value = "https://teddylee777.github.io/visualization/matplotlib-colorcode"

$ws.Range("D9").Value = "2022학년도 3월 신입/편입 1차 설명회"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/spring-2022-1st-webinar/#utm_source=rss&utm_medium=rss&utm_campaign=spring-2022-1st-webinar"

$ws.Range("D23").Value = "[pycharm파이참에서 import 자동으로 hide숨김을 해제하는 방법]How do I make pyCharm stop hiding (unfold) my Python imports?"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2924"

$ws.Range("D29").Value = "[만화] 인턴일기 58~65"
